# Prasad Commit @08:07 PM 08 August 2020
# Rewrites the header row of the "AddPayer" sheet to lowercase field
# names, widens column R, and updates the saved selection/scroll
# position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddPayer")

# --- Header row (row 1): same columns, new (mostly lower-cased) text ---
$ws.Range("A1").Value = "type"
$ws.Range("B1").Value = " ein"
$ws.Range("C1").Value = " ssn"
$ws.Range("D1").Value = "firstname"
$ws.Range("E1").Value = "middlename"
$ws.Range("F1").Value = " lastname"
$ws.Range("G1").Value = " businessname"
$ws.Range("H1").Value = " address"
$ws.Range("I1").Value = " suffix"
$ws.Range("J1").Value = " line2"
$ws.Range("K1").Value = " city"
$ws.Range("L1").Value = " state "
$ws.Range("M1").Value = "zipcode"
$ws.Range("N1").Value = " country"
$ws.Range("O1").Value = " checkheretoforeignaddress"
$ws.Range("P1").Value = " phone"
$ws.Range("Q1").Value = " email"
$ws.Range("R1").Value = " withholdingortaxstateId"
$ws.Range("S1").Value = " lastfiling"
$ws.Range("T1").Value = " clientid"

# --- Widen column R (18th column) ---
$ws.Columns("R").ColumnWidth = 29.5

# --- Move the view / selection to match the new scroll position ---
[void]$ws.Activate()
[void]$ws.Range("T3").Select()
$excel.ActiveWindow.ScrollColumn = 9   # column I -> matches topLeftCell="I1"
$excel.ActiveWindow.ScrollRow = 1
